# Updates to "Test Case" sheet:
#  - Remove the "True Results" column header (D2) and the whole "14) Add
#    favorite filter/identifier" test case row (row 17), which are no
#    longer tracked.
#  - Mark several more test-case rows as Pass ("x") in column E (and one
#    in column F), reflecting newly completed testing.
#  - Rename the "Update" test-case (row 9) and search-field expected
#    result (row 7) to reflect the new "Save"/"Note Search" behaviour.
#  - Leave the active selection on D2:D17 to match the reviewed range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the "True Results" column header - that column is no longer used.
$ws.Range("D2").Value = $null

# Row 4 (1 - UI functional) now passes.
$ws.Range("E4").Value = "x"

# Row 7 (4 - Search field test): expected result text updated, now passes.
$ws.Range("C7").Value = "Allows input of Note Search"
$ws.Range("E7").Value = "x"

# Row 8 (5 - Sample string search) now passes.
$ws.Range("E8").Value = "x"

# Row 9 (6 - renamed from Update Button test to Save Button test) now passes.
$ws.Range("B9").Value = "Save Button test"
$ws.Range("C9").Value = "Save results in text area to Note"
$ws.Range("E9").Value = "x"

# Row 12 (9 - Toggle button test) now passes (marked under Fail column).
$ws.Range("F12").Value = "x"

# Row 16 (13 - Display word count) now passes.
$ws.Range("E16").Value = "x"

# Remove the old "14) Add favorite filter/identifier" test case entirely.
$ws.Range("A17:C17").Value = $null

# Match the reviewed selection range.
$ws.Range("D2:D17").Select() | Out-Null
